# "add new line babui"
#
# Original last paragraph:
#   " efewfwefwef"   (single run)  + the _GoBack bookmark
#
# Target:
#   paragraph 1: " " / "E" / "fewfwefwef"   (same text, now split over 3 runs)
#   paragraph 2: empty paragraph
#   paragraph 3: "Hoisa babui"  (carries the _GoBack bookmark now)

$d = $word.ActiveDocument

# Locate the paragraph that currently holds " efewfwefwef" (the last
# paragraph in the body before the sectPr). Paragraph.Range.Text carries a
# trailing paragraph-mark character, so trim it before comparing.
$target = $d.Paragraphs.Last
if ($target.Range.Text.TrimEnd([char]13) -ne " efewfwefwef") {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13) -eq " efewfwefwef") {
            $target = $d.Paragraphs($i)
        }
    }
}

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Re-express " efewfwefwef" as three separate runs (" " / "E" / "fewfwefwef"),
# append a blank paragraph, then a new "Hoisa babui" paragraph that takes
# over the _GoBack bookmark. InsertXML lets us lay down exact run/paragraph
# boundaries instead of Word's usual same-format run merging.
$xml = "<w:p $w>" `
     +   "<w:r><w:t xml:space='preserve'> </w:t></w:r>" `
     +   "<w:r><w:t>E</w:t></w:r>" `
     +   "<w:r><w:t>fewfwefwef</w:t></w:r>" `
     + "</w:p>" `
     + "<w:p $w/>" `
     + "<w:p $w>" `
     +   "<w:r><w:t>Hoisa babui</w:t></w:r>" `
     +   "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" `
     +   "<w:bookmarkEnd w:id='0'/>" `
     + "</w:p>"

$target.Range.InsertXML($xml) | Out-Null
